$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 12: replace the old SW1/SW2/SW3 push-button BOM line with the
# new J3 pin-header connector line.
$ws.Range("A12").Value = "J3"
$ws.Range("D12").Value = "ZX-PZ2.54-2-4PWZ"
$ws.Range("C12").Value = "Connector_PinHeader_2.54mm:PinHeader_1x04_P2.54mm_Horizontal"
$ws.Range("B12").Value = "Conn_01x04"

# Apply the built-in "Hyperlink" look (underline, theme color 10) to the
# footprint cell without leaving an actual hyperlink relationship behind.
$ws.Hyperlinks.Add($ws.Range("C12"), "", "", "", "Connector_PinHeader_2.54mm:PinHeader_1x04_P2.54mm_Horizontal") | Out-Null
$ws.Hyperlinks.Delete()

# Update selection to match saved view state
$ws.Range("B12").Select()
